$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the timestamp on the existing last row (97) ---
$ws.Range("A97").Value = 45482.2916666667

# --- Append the new data row (98) reported by the R script ---
# Reuse row 97's date/time format for the new A98 cell (copy formats only,
# then overwrite with the new value) so no duplicate style gets minted.
$ws.Range("A97").Copy()
$ws.Range("A98").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A98").Value = 45483.6494675926

$ws.Range("B98").Value = 22500
$ws.Range("C98").Value = 3.3199999332428
$ws.Range("D98").Value = 3.00999999046326
$ws.Range("E98").Value = 3.19000005722046
$ws.Range("F98").Value = 3.33999991416931

# adj_close is stored as text matching the close price; force text so the
# numeric-looking value isn't auto-converted to a number, then drop the
# style back to the sheet's default (this column carries no explicit style).
$ws.Range("G98").Value = "'3.33999991416931"
$ws.Range("G98").Style = "Normal"

$ws.Range("H98").Value = "ESPE.MI"
